$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 149.6
$ws.Range("I8").Value = 149.6
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 448.8
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = -309.8
$ws.Range("N8").ClearContents()

$ws.Range("H62").Value = 4750.35
$ws.Range("I62").Value = 4378.875
$ws.Range("K62").Value = 4378.875
$ws.Range("M62").Value = -3754.875

$ws.Range("H65").Value = 4750.35
$ws.Range("I65").Value = 4378.875
$ws.Range("K65").Value = 21894.375
$ws.Range("M65").Value = -18774.375

$ws.Range("H116").Value = 7242.3
$ws.Range("J116").Value = 3939.3333
$ws.Range("L116").Value = 3939.3333
$ws.Range("N116").Value = -10823.3333

$ws.Range("H132").Value = 3379.5
$ws.Range("I132").Value = 1662.9722
$ws.Range("K132").Value = 4988.9166
$ws.Range("M132").Value = -2458.9166

$ws.Range("H135").Value = 203119.6
$ws.Range("I135").Value = 3899.75
$ws.Range("K135").Value = 35097.75
$ws.Range("M135").Value = -32562.75

$ws.Range("H138").Value = 2703.0645
$ws.Range("J138").Value = 2358.818
$ws.Range("L138").Value = 7076.454000000001
$ws.Range("N138").Value = -17356.454

$ws.Range("H141").Value = 8489.817999999999
$ws.Range("I141").Value = 7332.1113
$ws.Range("J141").Value = 13699.5
$ws.Range("K141").Value = 21996.3339
$ws.Range("L141").Value = 41098.5
$ws.Range("M141").Value = -16816.3339
$ws.Range("N141").Value = -51458.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H11").Value = 10475
$ws.Range("I11").Value = 950
$ws.Range("J11").Value = 20000
$ws.Range("K11").Value = 950
$ws.Range("L11").Value = 20000
$ws.Range("M11").Value = -806
$ws.Range("N11").Value = -20288

$ws.Range("H32").Value = 2683.2263
$ws.Range("I32").Value = 2683.2263
$ws.Range("K32").Value = 2683.2263
$ws.Range("M32").Value = -2396.2263

$ws.Range("H45").Value = 6549.3213
$ws.Range("I45").Value = 8788.375
$ws.Range("J45").Value = 3563.9167
$ws.Range("K45").Value = 8788.375
$ws.Range("L45").Value = 3563.9167
$ws.Range("M45").Value = -8411.375
$ws.Range("N45").Value = -4317.9167

$ws.Range("H132").Value = 1541.8214
$ws.Range("I132").Value = 1500.35
$ws.Range("K132").Value = 4501.049999999999
$ws.Range("M132").Value = -1971.049999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 9895.652
$ws.Range("I20").Value = 10716.353
$ws.Range("K20").Value = 10716.353
$ws.Range("M20").Value = -10469.353

$ws.Range("H99").Value = 2307.8333
$ws.Range("I99").Value = 2307.8333
$ws.Range("K99").Value = 2307.8333
$ws.Range("M99").Value = -809.8332999999998

$ws.Range("H105").Value = 6886.1904
$ws.Range("I105").Value = 10000.833
$ws.Range("K105").Value = 10000.833
$ws.Range("M105").Value = -8253.833000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H19").Value = 858.375
$ws.Range("I19").Value = 301.72726
$ws.Range("J19").Value = 2083
$ws.Range("K19").Value = 301.72726
$ws.Range("L19").Value = 2083
$ws.Range("M19").Value = -131.72726
$ws.Range("N19").Value = -2423

$ws.Range("H24").Value = 858.375
$ws.Range("I24").Value = 301.72726
$ws.Range("J24").Value = 2083
$ws.Range("K24").Value = 301.72726
$ws.Range("L24").Value = 2083
$ws.Range("M24").Value = -131.72726
$ws.Range("N24").Value = -2423

$ws.Range("H31").Value = 53573.31
$ws.Range("I31").Value = 56056.05
$ws.Range("K31").Value = 56056.05
$ws.Range("M31").Value = -55761.05

$ws.Range("H34").Value = 53573.31
$ws.Range("I34").Value = 56056.05
$ws.Range("K34").Value = 56056.05
$ws.Range("M34").Value = -55854.05

$ws.Range("H115").Value = 59156.332
$ws.Range("J115").Value = 59156.332
$ws.Range("L115").Value = 59156.332
$ws.Range("N115").Value = -61506.332

$ws.Range("H134").Value = 5455.273
$ws.Range("I134").Value = 4556.4443
$ws.Range("K134").Value = 13669.3329
$ws.Range("M134").Value = -11134.3329

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 113
$ws.Range("I12").Value = 64.2
$ws.Range("K12").Value = 192.6
$ws.Range("M12").Value = -19.60000000000002

$ws.Range("H25").Value = 200
$ws.Range("I25").Value = 200
$ws.Range("J25").Value = 0
$ws.Range("K25").Value = 600
$ws.Range("L25").Value = 0
$ws.Range("M25").Value = -431
$ws.Range("N25").ClearContents()

$ws.Range("H30").Value = 200
$ws.Range("I30").Value = 200
$ws.Range("J30").Value = 0
$ws.Range("K30").Value = 600
$ws.Range("L30").Value = 0
$ws.Range("M30").Value = -498
$ws.Range("N30").ClearContents()

$ws.Range("H60").Value = 788
$ws.Range("I60").Value = 482.85715
$ws.Range("J60").Value = 1500
$ws.Range("K60").Value = 1448.57145
$ws.Range("L60").Value = 4500
$ws.Range("M60").Value = -1197.57145
$ws.Range("N60").Value = -5002

$ws.Range("H132").Value = 2530.35
$ws.Range("I132").Value = 1899.8334
$ws.Range("J132").Value = 2800.5715
$ws.Range("K132").Value = 17098.5006
$ws.Range("L132").Value = 25205.1435
$ws.Range("M132").Value = -14568.5006
$ws.Range("N132").Value = -30265.1435

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6776.4443
$ws.Range("I70").Value = 6081.3335
$ws.Range("K70").Value = 6081.3335
$ws.Range("M70").Value = -5811.3335

$ws.Range("H73").Value = 6776.4443
$ws.Range("I73").Value = 6081.3335
$ws.Range("K73").Value = 6081.3335
$ws.Range("M73").Value = -5145.3335

$ws.Range("H107").Value = 21743020
$ws.Range("I107").Value = 1814.9
$ws.Range("K107").Value = 1814.9
$ws.Range("M107").Value = 105.0999999999999

$ws.Range("H113").Value = 4334.5454
$ws.Range("J113").Value = 3477
$ws.Range("L113").Value = 3477
$ws.Range("N113").Value = -7817

$ws.Range("H132").Value = 4394.102
$ws.Range("I132").Value = 3620.8604
$ws.Range("K132").Value = 10862.5812
$ws.Range("M132").Value = -8332.581200000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2395.3635
$ws.Range("I46").Value = 2500
$ws.Range("J46").Value = 2372.111
$ws.Range("K46").Value = 2500
$ws.Range("L46").Value = 2372.111
$ws.Range("M46").Value = -2312
$ws.Range("N46").Value = -2748.111

$ws.Range("H132").Value = 3599.08
$ws.Range("I132").Value = 3162.818
$ws.Range("K132").Value = 9488.454000000002
$ws.Range("M132").Value = -6958.454000000002

$ws.Range("H136").Value = 5285.615
$ws.Range("I136").Value = 4883.227
$ws.Range("K136").Value = 14649.681
$ws.Range("M136").Value = -12099.681

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H31").Value = 0
$ws.Range("I31").Value = 0
$ws.Range("K31").Value = 0
$ws.Range("M31").ClearContents()

$ws.Range("H88").Value = 24517.75
$ws.Range("I88").Value = 34085.5
$ws.Range("J88").Value = 14950
$ws.Range("K88").Value = 34085.5
$ws.Range("L88").Value = 14950
$ws.Range("M88").Value = -33679.5
$ws.Range("N88").Value = -15762

$ws.Range("H91").Value = 24517.75
$ws.Range("I91").Value = 34085.5
$ws.Range("J91").Value = 14950
$ws.Range("K91").Value = 34085.5
$ws.Range("L91").Value = 14950
$ws.Range("M91").Value = -32681.5
$ws.Range("N91").Value = -17758

$ws.Range("H132").Value = 11086.286
$ws.Range("I132").Value = 15800.889
$ws.Range("K132").Value = 47402.667
$ws.Range("M132").Value = -44872.667
